$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Repurpose the existing "总计" sheet (currently sheet index 3) as
#    the new "2022-Q1" detail sheet, keeping its sheetId/r:id so the
#    new "总计" sheet can be appended at the end with a fresh id.
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"

# Propagate the existing header style (currently on B1:D1) onto the
# new header cells E1:H1 and the existing index-column style (on
# A2:A3) onto the new A4:A5 cells, via copy/paste-of-formats so we
# reuse the same style index instead of inventing new ones.
$q1.Range("B1").Copy() | Out-Null
$q1.Range("E1:H1").PasteSpecial(-4122) | Out-Null
$q1.Range("A2").Copy() | Out-Null
$q1.Range("A4:A5").PasteSpecial(-4122) | Out-Null
$q1.Application.CutCopyMode = $false

# ---- Header row ----
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Columns B and D:G hold numeric-looking figures (fund codes, fund
# size, position weight, ...) that must be stored as TEXT (matching
# the workbook's existing convention on every other detail sheet, and
# preserving leading zeros in fund codes). Force text via
# NumberFormat, write the values, then drop back to the default
# "Normal" style so no stray s="" survives on the cell (NumberFormat
# reverts to General / style 0). Column C (fund name) is included too
# since it sits inside the bounding box, but it is never numeric-
# looking so this is a no-op for it.
$textCells = $q1.Range("B2:G5")
$textCells.NumberFormat = "@"

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "011052"
$q1.Range("C2").Value = "鹏华弘裕一年持有期混合A"
$q1.Range("D2").Value = "2.92"
$q1.Range("E2").Value = "24.56"
$q1.Range("F2").Value = "3.42"
$q1.Range("G2").Value = "0.0999"
$q1.Range("H2").Value = 1

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "001190"
$q1.Range("C3").Value = "鹏华弘润灵活配置混合 - A"
$q1.Range("D3").Value = "3.65"
$q1.Range("E3").Value = "23.00"
$q1.Range("F3").Value = "1.48"
$q1.Range("G3").Value = "0.0540"
$q1.Range("H3").Value = 3

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "011053"
$q1.Range("C4").Value = "鹏华弘裕一年持有期混合C"
$q1.Range("D4").Value = "0.20"
$q1.Range("E4").Value = "24.56"
$q1.Range("F4").Value = "3.42"
$q1.Range("G4").Value = "0.0068"
$q1.Range("H4").Value = 1

$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "001191"
$q1.Range("C5").Value = "鹏华弘润灵活配置混合 - C"
$q1.Range("D5").Value = "0.25"
$q1.Range("E5").Value = "23.00"
$q1.Range("F5").Value = "1.48"
$q1.Range("G5").Value = "0.0037"
$q1.Range("H5").Value = 3

$textCells.Style = "Normal"

# ------------------------------------------------------------------
# 2. Append a brand-new "总计" summary sheet at the end of the book.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add($null, $lastSheet)
$total.Name = "总计"

# Reuse the same header / index-column style (s="2") that already
# lives in the style table (carried by the "2022-Q1" sheet we just
# repurposed) instead of building the formatting up by hand, so no
# new style entries are introduced.
$q1.Range("B1:D1").Copy() | Out-Null
$total.Range("B1:D1").PasteSpecial(-4122) | Out-Null
$q1.Range("A2").Copy() | Out-Null
$total.Range("A2:A4").PasteSpecial(-4122) | Out-Null
$total.Application.CutCopyMode = $false

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.16

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 6
$total.Range("D3").Value = 0.27

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q2"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.18
